# cv_burial.xlsx - add missing "Burial Filling Objects" / "Burial Filling Type"
# autocomplete rows to the concept-value sheet (rows 207-217).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is (A = conceptScheme / header label, B = pref_label@eng | pref_label@ger)
$newRows = @(
    @("Burial Filling Objects| Grabfüllungsobjekte", "pebble stone|Kieselstein"),
    @("Burial Filling Objects| Grabfüllungsobjekte", "bigger stone|größere Steine"),
    @("Burial Filling Objects| Grabfüllungsobjekte", "ceramic sherd|Scherbenschüttung (Keramikbruchstücke?)"),
    @("Burial Filling Objects| Grabfüllungsobjekte", "secondary cremated ceramic sherds|Scherben - sekundär gebrannt"),
    @("Burial Filling Objects| Grabfüllungsobjekte", "human remains|Menschenknochen"),
    @("Burial Filling Objects| Grabfüllungsobjekte", "animal remains|Tierknochen"),
    @("Burial Filling Objects| Grabfüllungsobjekte", "snail/shell|Schnecke/Muschel"),
    @("Burial Filling Type|Art der Grabfüllung", "monophase and unburnt|einphasig – ungebrannt"),
    @("Burial Filling Type|Art der Grabfüllung", "monophase and burnt|einphasig – gebrannt"),
    @("Burial Filling Type|Art der Grabfüllung", "multi phase and unburnt|mehrphasig – ungebrannt"),
    @("Burial Filling Type|Art der Grabfüllung", "multi phase and burnt|mehrphasig – gebrannt")
)

$startRow = 207
$r = $startRow
foreach ($pair in $newRows) {
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
    $ws.Cells.Item($r, 2).WrapText = $true
    $r = $r + 1
}

$lastRow = $r - 1

# Keep the visible selection/scroll near the newly-appended rows, mirroring
# the author's last edit position (bottom of the new block).
$ws.Range("B" + $lastRow).Select()

Write-Host ("Added rows " + $startRow + ":" + $lastRow)
